# "Generate Report for Handoff"
#
# Updates the localization-status report for the "b.md" file now that a new
# handoff (b.*.xlf) has been generated for it: status moves from
# "Handed back: in sync with en-US" to "Ready for handoff", the latest
# handoff file/datetime are refreshed, content-duplicate flips to False, and
# an error detail message is recorded because the existing handback is stale.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) - Status columns (zh-cn/de-de) and the
# "Latest HO Xliff Generate Date" column.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 07:15:46"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# NOTE: assigning the bare text "False" gets auto-coerced to a native
# boolean by the Value setter. Prefix with an apostrophe to force text,
# then restore the (unstyled) Normal cell style so formatting matches the
# rest of the sheet.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-01-03 07:15:35"
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a4b80ca38a28cdfb7ac57ae17c50f2577c2c14dc/e2e/b.md."
$wsZhCn.Columns.Item(18).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-01-03 07:15:46"
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a2a32c0a5631a13868300f9802c4f19155acb3b9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test1/blob/a4b80ca38a28cdfb7ac57ae17c50f2577c2c14dc/e2e/b.md."
$wsDeDe.Columns.Item(18).ColumnWidth = 39.15
